$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 69

$ws.Cells.Item($row, 1).Value = "XNFKAQ"
$ws.Cells.Item($row, 2).Value = "Film de fusor HP"
$ws.Cells.Item($row, 3).Value = "1600 2600 2605, CP1025 CP1215 CP1515 CP1518 CP1525 CP2025, M175 M176 M177 M275 M276 M351 M375 M451 M475 M476, CM1015 CM1017 CM1415 CM1312 CM2320"
$ws.Cells.Item($row, 4).Value = 50000
$ws.Cells.Item($row, 5).Value = 150000
$ws.Cells.Item($row, 6).Value = 9
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Formula = "=(E69-D69)*G69"
$ws.Cells.Item($row, 9).Formula = "=D69*F69"
$ws.Cells.Item($row, 10).Value = 450000
